$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C26").Value = "V"
$ws.Range("E26").Value = "V"
$ws.Range("F26").Value = "V"
$ws.Range("G26").Value = "V"
$ws.Range("H26").Value = "X"
$ws.Range("K26").Value = "V"
